# Canopy Aerospace timesheet (Calvin) - add 4 new weekly blocks (weeks of
# 11/12, 11/19, 11/26, 12/3) below the existing data, following the same
# template as the prior weeks, and fix up a couple of numbers in the most
# recent existing week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing week (rows 102:113): Sponsor Meeting gained 1 hour on
#     Wednesday (E106), and the week's Wednesday / Week-Total columns
#     change accordingly (E113, I113).
$ws.Range("E106").Value = 1
$ws.Range("I106").Value = 1
$ws.Range("E113").Value = 5
$ws.Range("I113").Value = 17

# --- Add four new weekly blocks by copying the prior week's 13-row
#     template (header band + date row + 9 task rows + daily-total row)
#     down to each new location. Using Range.Copy(Destination) (single
#     call form) keeps styles/number-formats/merge intact.
$ws.Range("A101:I113").Copy($ws.Range("A116"))
$ws.Range("A101:I113").Copy($ws.Range("A131"))
$ws.Range("A101:I113").Copy($ws.Range("A146"))
$ws.Range("A101:I113").Copy($ws.Range("A161"))

# --- Week of 11/12 - 11/18 (rows 116-128) ---
$ws.Range("B117").Value = 45608
$ws.Range("C117").Value = 45609
$ws.Range("D117").Value = 45610
$ws.Range("E117").Value = 45245
$ws.Range("F117").Value = 45246
$ws.Range("G117").Value = 45247
$ws.Range("H117").Value = 45248

$ws.Range("A118").Value = "Lecture"
$ws.Range("B118").Value = ""
$ws.Range("C118").Value = ""
$ws.Range("D118").Value = ""
$ws.Range("E118").Value = ""
$ws.Range("F118").Value = ""
$ws.Range("G118").Value = ""
$ws.Range("I118").Value = 0

$ws.Range("A119").Value = "Read/Study"
$ws.Range("B119").Value = ""
$ws.Range("C119").Value = ""
$ws.Range("D119").Value = ""
$ws.Range("E119").Value = ""
$ws.Range("F119").Value = ""
$ws.Range("G119").Value = ""
$ws.Range("I119").Value = 0

$ws.Range("A120").Value = "Team Meeting"
$ws.Range("B120").Value = ""
$ws.Range("C120").Value = 1
$ws.Range("D120").Value = ""
$ws.Range("E120").Value = ""
$ws.Range("F120").Value = ""
$ws.Range("G120").Value = ""
$ws.Range("I120").Value = 1

$ws.Range("A121").Value = "Sponsor Meeting"
$ws.Range("B121").Value = ""
$ws.Range("C121").Value = ""
$ws.Range("D121").Value = ""
$ws.Range("E121").Value = 1
$ws.Range("F121").Value = ""
$ws.Range("G121").Value = ""
$ws.Range("I121").Value = 1

$ws.Range("A122").Value = "Setup dev environment"
$ws.Range("B122").Value = ""
$ws.Range("C122").Value = ""
$ws.Range("D122").Value = ""
$ws.Range("E122").Value = ""
$ws.Range("F122").Value = ""
$ws.Range("G122").Value = ""
$ws.Range("I122").Value = 0

$ws.Range("A123").Value = "Read literature"
$ws.Range("B123").Value = ""
$ws.Range("C123").Value = ""
$ws.Range("D123").Value = ""
$ws.Range("E123").Value = ""
$ws.Range("F123").Value = ""
$ws.Range("G123").Value = ""
$ws.Range("I123").Value = 0

$ws.Range("A124").Value = "Admin tasks"
$ws.Range("B124").Value = ""
$ws.Range("C124").Value = ""
$ws.Range("D124").Value = ""
$ws.Range("E124").Value = ""
$ws.Range("F124").Value = ""
$ws.Range("G124").Value = ""
$ws.Range("I124").Value = 0

$ws.Range("A125").Value = "Visualization"
$ws.Range("B125").Value = ""
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 1
$ws.Range("E125").Value = 2
$ws.Range("F125").Value = 1
$ws.Range("G125").Value = 1
$ws.Range("I125").Value = 6

$ws.Range("A126").Value = "Task yyyyy"
$ws.Range("B126").Value = ""
$ws.Range("C126").Value = ""
$ws.Range("D126").Value = ""
$ws.Range("E126").Value = ""
$ws.Range("F126").Value = ""
$ws.Range("G126").Value = ""
$ws.Range("I126").Value = 0

$ws.Range("A127").Value = "Task zzzzzz"
$ws.Range("B127").Value = ""
$ws.Range("C127").Value = ""
$ws.Range("D127").Value = ""
$ws.Range("E127").Value = ""
$ws.Range("F127").Value = ""
$ws.Range("G127").Value = ""
$ws.Range("I127").Value = 0

$ws.Range("B128").Value = 0
$ws.Range("C128").Value = 2
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = 3
$ws.Range("F128").Value = 1
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 8

# --- Week of 11/19 - 11/25 (rows 131-143) - Thanksgiving week, mostly
#     empty except a "THANKSGIVING!!!" note on Wednesday ---
$ws.Range("B132").Value = 45615
$ws.Range("C132").Value = 45616
$ws.Range("D132").Value = 45617
$ws.Range("E132").Value = 45252
$ws.Range("F132").Value = 45253
$ws.Range("G132").Value = 45254
$ws.Range("H132").Value = 45255

$ws.Range("A133").Value = "Lecture"
$ws.Range("B133").Value = ""
$ws.Range("C133").Value = ""
$ws.Range("D133").Value = ""
$ws.Range("E133").Value = ""
$ws.Range("F133").Value = ""
$ws.Range("G133").Value = ""
$ws.Range("I133").Value = 0

$ws.Range("A134").Value = "Read/Study"
$ws.Range("B134").Value = ""
$ws.Range("C134").Value = ""
$ws.Range("D134").Value = ""
$ws.Range("E134").Value = ""
$ws.Range("F134").Value = ""
$ws.Range("G134").Value = ""
$ws.Range("I134").Value = 0

$ws.Range("A135").Value = "Team Meeting"
$ws.Range("B135").Value = ""
$ws.Range("C135").Value = ""
$ws.Range("D135").Value = ""
$ws.Range("E135").Value = ""
$ws.Range("F135").Value = ""
$ws.Range("G135").Value = ""
$ws.Range("I135").Value = 0

$ws.Range("A136").Value = "Sponsor Meeting"
$ws.Range("B136").Value = ""
$ws.Range("C136").Value = ""
$ws.Range("D136").Value = ""
$ws.Range("E136").Value = ""
$ws.Range("F136").Value = ""
$ws.Range("G136").Value = ""
$ws.Range("I136").Value = 0

$ws.Range("A137").Value = "Setup dev environment"
$ws.Range("B137").Value = ""
$ws.Range("C137").Value = ""
$ws.Range("D137").Value = "THANKSGIVING!!!"
$ws.Range("E137").Value = ""
$ws.Range("F137").Value = ""
$ws.Range("G137").Value = ""
$ws.Range("I137").Value = 0

$ws.Range("A138").Value = "Read literature"
$ws.Range("B138").Value = ""
$ws.Range("C138").Value = ""
$ws.Range("D138").Value = ""
$ws.Range("E138").Value = ""
$ws.Range("F138").Value = ""
$ws.Range("G138").Value = ""
$ws.Range("I138").Value = 0

$ws.Range("A139").Value = "Admin tasks"
$ws.Range("B139").Value = ""
$ws.Range("C139").Value = ""
$ws.Range("D139").Value = ""
$ws.Range("E139").Value = ""
$ws.Range("F139").Value = ""
$ws.Range("G139").Value = ""
$ws.Range("I139").Value = 0

$ws.Range("A140").Value = "Visualization"
$ws.Range("B140").Value = ""
$ws.Range("C140").Value = ""
$ws.Range("D140").Value = ""
$ws.Range("E140").Value = ""
$ws.Range("F140").Value = ""
$ws.Range("G140").Value = ""
$ws.Range("I140").Value = 0

$ws.Range("A141").Value = "Task yyyyy"
$ws.Range("B141").Value = ""
$ws.Range("C141").Value = ""
$ws.Range("D141").Value = ""
$ws.Range("E141").Value = ""
$ws.Range("F141").Value = ""
$ws.Range("G141").Value = ""
$ws.Range("I141").Value = 0

$ws.Range("A142").Value = "Task zzzzzz"
$ws.Range("B142").Value = ""
$ws.Range("C142").Value = ""
$ws.Range("D142").Value = ""
$ws.Range("E142").Value = ""
$ws.Range("F142").Value = ""
$ws.Range("G142").Value = ""
$ws.Range("I142").Value = 0

$ws.Range("B143").Value = 0
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 0
$ws.Range("E143").Value = 0
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 0
$ws.Range("I143").Value = 0

# --- Week of 11/26 - 12/2 (rows 146-158) ---
$ws.Range("B147").Value = 45622
$ws.Range("C147").Value = 45623
$ws.Range("D147").Value = 45624
$ws.Range("E147").Value = 45259
$ws.Range("F147").Value = 45260
$ws.Range("G147").Value = 45261
$ws.Range("H147").Value = 45262

$ws.Range("A148").Value = "Lecture"
$ws.Range("B148").Value = ""
$ws.Range("C148").Value = ""
$ws.Range("D148").Value = ""
$ws.Range("E148").Value = ""
$ws.Range("F148").Value = ""
$ws.Range("G148").Value = ""
$ws.Range("I148").Value = 0

$ws.Range("A149").Value = "Read/Study"
$ws.Range("B149").Value = ""
$ws.Range("C149").Value = ""
$ws.Range("D149").Value = ""
$ws.Range("E149").Value = ""
$ws.Range("F149").Value = ""
$ws.Range("G149").Value = ""
$ws.Range("I149").Value = 0

$ws.Range("A150").Value = "Team Meeting"
$ws.Range("B150").Value = ""
$ws.Range("C150").Value = 1
$ws.Range("D150").Value = ""
$ws.Range("E150").Value = ""
$ws.Range("F150").Value = ""
$ws.Range("G150").Value = ""
$ws.Range("I150").Value = 1

$ws.Range("A151").Value = "Sponsor Meeting"
$ws.Range("B151").Value = ""
$ws.Range("C151").Value = ""
$ws.Range("D151").Value = ""
$ws.Range("E151").Value = 1
$ws.Range("F151").Value = ""
$ws.Range("G151").Value = ""
$ws.Range("I151").Value = 1

$ws.Range("A152").Value = "Setup dev environment"
$ws.Range("B152").Value = ""
$ws.Range("C152").Value = ""
$ws.Range("D152").Value = ""
$ws.Range("E152").Value = ""
$ws.Range("F152").Value = ""
$ws.Range("G152").Value = ""
$ws.Range("I152").Value = 0

$ws.Range("A153").Value = "Read literature"
$ws.Range("B153").Value = ""
$ws.Range("C153").Value = ""
$ws.Range("D153").Value = ""
$ws.Range("E153").Value = ""
$ws.Range("F153").Value = ""
$ws.Range("G153").Value = ""
$ws.Range("I153").Value = 0

$ws.Range("A154").Value = "Admin tasks"
$ws.Range("B154").Value = ""
$ws.Range("C154").Value = ""
$ws.Range("D154").Value = ""
$ws.Range("E154").Value = 1
$ws.Range("F154").Value = ""
$ws.Range("G154").Value = ""
$ws.Range("I154").Value = 1

$ws.Range("A155").Value = "Visualization"
$ws.Range("B155").Value = ""
$ws.Range("C155").Value = ""
$ws.Range("D155").Value = 2
$ws.Range("E155").Value = 1
$ws.Range("F155").Value = 3
$ws.Range("G155").Value = 1
$ws.Range("I155").Value = 7

$ws.Range("A156").Value = "Task yyyyy"
$ws.Range("B156").Value = ""
$ws.Range("C156").Value = ""
$ws.Range("D156").Value = ""
$ws.Range("E156").Value = ""
$ws.Range("F156").Value = ""
$ws.Range("G156").Value = ""
$ws.Range("I156").Value = 0

$ws.Range("A157").Value = "Task zzzzzz"
$ws.Range("B157").Value = ""
$ws.Range("C157").Value = ""
$ws.Range("D157").Value = ""
$ws.Range("E157").Value = ""
$ws.Range("F157").Value = ""
$ws.Range("G157").Value = ""
$ws.Range("I157").Value = 0

$ws.Range("B158").Value = 0
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 2
$ws.Range("E158").Value = 3
$ws.Range("F158").Value = 3
$ws.Range("G158").Value = 1
$ws.Range("H158").Value = 0
$ws.Range("I158").Value = 10

# --- Week of 12/3 - 12/9 (rows 161-173) ---
$ws.Range("B162").Value = 45629
$ws.Range("C162").Value = 45630
$ws.Range("D162").Value = 45631
$ws.Range("E162").Value = 45266
$ws.Range("F162").Value = 45267
$ws.Range("G162").Value = 45268
$ws.Range("H162").Value = 45269

$ws.Range("A163").Value = "Lecture"
$ws.Range("B163").Value = ""
$ws.Range("C163").Value = ""
$ws.Range("D163").Value = ""
$ws.Range("E163").Value = ""
$ws.Range("F163").Value = ""
$ws.Range("G163").Value = ""
$ws.Range("I163").Value = 0

$ws.Range("A164").Value = "Read/Study"
$ws.Range("B164").Value = ""
$ws.Range("C164").Value = ""
$ws.Range("D164").Value = ""
$ws.Range("E164").Value = ""
$ws.Range("F164").Value = ""
$ws.Range("G164").Value = ""
$ws.Range("I164").Value = 0

$ws.Range("A165").Value = "Team Meeting"
$ws.Range("B165").Value = ""
$ws.Range("C165").Value = 1
$ws.Range("D165").Value = ""
$ws.Range("E165").Value = ""
$ws.Range("F165").Value = ""
$ws.Range("G165").Value = ""
$ws.Range("I165").Value = 1

$ws.Range("A166").Value = "Sponsor Meeting"
$ws.Range("B166").Value = ""
$ws.Range("C166").Value = ""
$ws.Range("D166").Value = ""
$ws.Range("E166").Value = ""
$ws.Range("F166").Value = ""
$ws.Range("G166").Value = ""
$ws.Range("I166").Value = 0

$ws.Range("A167").Value = "Setup dev environment"
$ws.Range("B167").Value = ""
$ws.Range("C167").Value = ""
$ws.Range("D167").Value = ""
$ws.Range("E167").Value = ""
$ws.Range("F167").Value = ""
$ws.Range("G167").Value = ""
$ws.Range("I167").Value = 0

$ws.Range("A168").Value = "Read literature"
$ws.Range("B168").Value = ""
$ws.Range("C168").Value = ""
$ws.Range("D168").Value = ""
$ws.Range("E168").Value = ""
$ws.Range("F168").Value = ""
$ws.Range("G168").Value = ""
$ws.Range("I168").Value = 0

$ws.Range("A169").Value = "Admin tasks"
$ws.Range("B169").Value = ""
$ws.Range("C169").Value = ""
$ws.Range("D169").Value = ""
$ws.Range("E169").Value = ""
$ws.Range("F169").Value = ""
$ws.Range("G169").Value = ""
$ws.Range("I169").Value = 0

$ws.Range("A170").Value = "Visualization"
$ws.Range("B170").Value = ""
$ws.Range("C170").Value = ""
$ws.Range("D170").Value = ""
$ws.Range("E170").Value = ""
$ws.Range("F170").Value = ""
$ws.Range("G170").Value = ""
$ws.Range("I170").Value = 0

$ws.Range("A171").Value = "Task yyyyy"
$ws.Range("B171").Value = ""
$ws.Range("C171").Value = ""
$ws.Range("D171").Value = ""
$ws.Range("E171").Value = ""
$ws.Range("F171").Value = ""
$ws.Range("G171").Value = ""
$ws.Range("I171").Value = 0

$ws.Range("A172").Value = "Task zzzzzz"
$ws.Range("B172").Value = ""
$ws.Range("C172").Value = ""
$ws.Range("D172").Value = ""
$ws.Range("E172").Value = ""
$ws.Range("F172").Value = ""
$ws.Range("G172").Value = ""
$ws.Range("I172").Value = 0

$ws.Range("B173").Value = 0
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 0
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 16

# --- Update the view so the newly added week is what's on screen, same
#     as the author left it after editing (matches the saved sheetView). ---
$ws.Range("G158").Select()
$excel.ActiveWindow.ScrollRow = 150
